{"js": "// Update the date heading and every multiplication problem in the\n// practice-sheet table to the values for the new day.\nconst replacements = [\n  [\"2026-02-03 Tuesday\", \"2026-02-04 Wednesday\"],\n  [\"93\u00d781=\", \"81\u00d748=\"],\n  [\"20\u00d795=\", \"48\u00d714=\"],\n  [\"53\u00d717=\", \"98\u00d727=\"],\n  [\"63\u00d754=\", \"21\u00d772=\"],\n  [\"52\u00d711=\", \"53\u00d767=\"],\n  [\"13\u00d753=\", \"28\u00d770=\"],\n  [\"56\u00d739=\", \"56\u00d750=\"],\n  [\"12\u00d770=\", \"63\u00d763=\"],\n  [\"99\u00d725=\", \"35\u00d789=\"],\n  [\"34\u00d788=\", \"27\u00d771=\"],\n  [\"39\u00d728=\", \"18\u00d741=\"],\n  [\"37\u00d713=\", \"54\u00d738=\"],\n  [\"38\u00d799=\", \"26\u00d736=\"],\n  [\"58\u00d754=\", \"79\u00d724=\"],\n  [\"35\u00d795=\", \"48\u00d755=\"],\n  [\"57\u00d745=\", \"55\u00d794=\"],\n  [\"62\u00d799=\", \"75\u00d764=\"],\n  [\"99\u00d743=\", \"57\u00d725=\"],\n  [\"39\u00d712=\", \"89\u00d786=\"],\n  [\"85\u00d765=\", \"31\u00d747=\"],\n  [\"68\u00d747=\", \"68\u00d719=\"],\n  [\"63\u00d746=\", \"74\u00d788=\"],\n  [\"71\u00d718=\", \"33\u00d738=\"],\n  [\"58\u00d778=\", \"17\u00d747=\"],\n  [\"85\u00d758=\", \"94\u00d784=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Updates the date heading and all multiplication problems in the table\n# to match the new day's values, via Word's Find/Replace (wdReplaceAll).\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\nReplace-AllText '2026-02-03 Tuesday' '2026-02-04 Wednesday'\nReplace-AllText '93\u00d781=' '81\u00d748='\nReplace-AllText '20\u00d795=' '48\u00d714='\nReplace-AllText '53\u00d717=' '98\u00d727='\nReplace-AllText '63\u00d754=' '21\u00d772='\nReplace-AllText '52\u00d711=' '53\u00d767='\nReplace-AllText '13\u00d753=' '28\u00d770='\nReplace-AllText '56\u00d739=' '56\u00d750='\nReplace-AllText '12\u00d770=' '63\u00d763='\nReplace-AllText '99\u00d725=' '35\u00d789='\nReplace-AllText '34\u00d788=' '27\u00d771='\nReplace-AllText '39\u00d728=' '18\u00d741='\nReplace-AllText '37\u00d713=' '54\u00d738='\nReplace-AllText '38\u00d799=' '26\u00d736='\nReplace-AllText '58\u00d754=' '79\u00d724='\nReplace-AllText '35\u00d795=' '48\u00d755='\nReplace-AllText '57\u00d745=' '55\u00d794='\nReplace-AllText '62\u00d799=' '75\u00d764='\nReplace-AllText '99\u00d743=' '57\u00d725='\nReplace-AllText '39\u00d712=' '89\u00d786='\nReplace-AllText '85\u00d765=' '31\u00d747='\nReplace-AllText '68\u00d747=' '68\u00d719='\nReplace-AllText '63\u00d746=' '74\u00d788='\nReplace-AllText '71\u00d718=' '33\u00d738='\nReplace-AllText '58\u00d778=' '17\u00d747='\nReplace-AllText '85\u00d758=' '94\u00d784='\n"}
